$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Afganistan inserted right after Tailandia (before Grecia) ---
# Row 65: now Afganistan, with new data values
$ws.Cells.Item(65, 1).Value = "Afganistan"
$ws.Cells.Item(65, 2).Value = 2704
$ws.Cells.Item(65, 3).Value = 235
$ws.Cells.Item(65, 4).Value = 345
$ws.Cells.Item(65, 5).Value = 2274
$ws.Cells.Item(65, 6).Value = 7
$ws.Cells.Item(65, 7).Value = 13
$ws.Cells.Item(65, 8).Value = 85

# Row 66: now Grecia, shifted down from old row 65
$ws.Cells.Item(66, 1).Value = "Grecia"
$ws.Cells.Item(66, 2).Value = 2620
$ws.Cells.Item(66, 3).Value = 0
$ws.Cells.Item(66, 4).Value = 1374
$ws.Cells.Item(66, 5).Value = 1103
$ws.Cells.Item(66, 6).Value = 37
$ws.Cells.Item(66, 7).Value = 0
$ws.Cells.Item(66, 8).Value = 143

# Row 67: now Oman, shifted down from old row 66
$ws.Cells.Item(67, 1).Value = "Oman"
$ws.Cells.Item(67, 2).Value = 2568
$ws.Cells.Item(67, 3).Value = 85
$ws.Cells.Item(67, 4).Value = 750
$ws.Cells.Item(67, 5).Value = 1806
$ws.Cells.Item(67, 6).Value = 17
$ws.Cells.Item(67, 7).Value = 0
$ws.Cells.Item(67, 8).Value = 12

# Row 68 (Nigeria) is unchanged.

# --- Albania inserted right after Republica de Chipre (before Kirguistan) ---
# Row 96: now Albania, with new data values
$ws.Cells.Item(96, 1).Value = "Albania"
$ws.Cells.Item(96, 2).Value = 795
$ws.Cells.Item(96, 3).Value = 6
$ws.Cells.Item(96, 4).Value = 531
$ws.Cells.Item(96, 5).Value = 233
$ws.Cells.Item(96, 6).Value = 7
$ws.Cells.Item(96, 7).Value = 0
$ws.Cells.Item(96, 8).Value = 31

# Row 97: now Kirguistan, shifted down from old row 96
$ws.Cells.Item(97, 1).Value = "Kirguistan"
$ws.Cells.Item(97, 2).Value = 795
$ws.Cells.Item(97, 3).Value = 26
$ws.Cells.Item(97, 4).Value = 564
$ws.Cells.Item(97, 5).Value = 221
$ws.Cells.Item(97, 6).Value = 12
$ws.Cells.Item(97, 7).Value = 2
$ws.Cells.Item(97, 8).Value = 10

# Row 98 (Principado de Andorra) is unchanged.
